$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137; this shifts existing rows 137:219 down to 138:220
# (Excel.Rows.Insert also copies formatting from the row above, matching s="2" on column D).
$ws.Rows("137").Insert()

# Populate the newly inserted row 137 with the new price-observation record.
$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = 44596
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 100112043
$ws.Range("G137").Value = "Pepino ensalada"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 400
$ws.Range("K137").Value = 17000
$ws.Range("L137").Value = 17000
$ws.Range("M137").Value = 17000
$ws.Range("N137").Value = "$/caja 60 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 283
$ws.Range("Q137").Value = 60
$ws.Range("R137").Value = "Hortaliza"
